# Update "Percent Complete" values in the completeness report.
# These cells hold text-typed numeric strings (e.g. "99.1"); we force
# text entry via a leading apostrophe (quote-prefix) and then restore the
# cell's style to "Normal" so no explicit style index is left behind,
# matching the original (unstyled) cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    $rng = $ws.Range($Cell)
    $rng.Value = "'" + $Text
    $rng.Style = "Normal"
}

# Column B ("Fields of Interest" block, rows 2-11)
Set-TextValue "B5"  "99.6"
Set-TextValue "B6"  "98.7"
Set-TextValue "B7"  "99.6"
Set-TextValue "B8"  "99.6"
Set-TextValue "B9"  "98.6"
Set-TextValue "B10" "98.3"
Set-TextValue "B11" "95.7"

# Column E ("Fields of Interest" block, rows 2-26)
Set-TextValue "E13" "96.9"
Set-TextValue "E15" "65.7"
Set-TextValue "E16" "62.0"
Set-TextValue "E17" "59.8"
Set-TextValue "E18" "58.3"
Set-TextValue "E19" "59.8"
Set-TextValue "E20" "61.4"

Write-Host "Completeness report percentages updated."
